# Weekly CompStat data refresh (70th Precinct): roll the report forward one
# week and update the crime-complaint figures with the newly collected data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header: bump "Volume 31  Number 27" -> "Number 28" and the reporting week
# from 7/1/2024-7/7/2024 to 7/8/2024-7/14/2024.
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 31   Number  28"
$ws.Range("C9").Value = "Report Covering the Week  7/8/2024  Through  7/14/2024"

# ---------------------------------------------------------------------------
# Helpers for the data table (rows 14-31).
#
# Most cells are plain numbers and can simply be overwritten in place
# (Set-Num). A handful of cells flip between a numeric value and the
# sheet's textual placeholders "0" / "***.*" (used when a rate can't be
# computed, e.g. dividing by a zero prior-period count). Those placeholders
# are shared strings already present in the workbook (index 20 = "0",
# index 21 = "***.*"), and the cells that hold them use the same "General"
# right-aligned style (s=14) as the rest of the text cells in the table, so
# we clone that style from an existing text cell before writing the value.
# ---------------------------------------------------------------------------

function Set-Num($cellRef, $val) {
    $ws.Range($cellRef).Value = $val
}

function Set-NumFrom($cellRef, $donorRef, $val) {
    # Re-use the numeric-cell style (copied from $donorRef) then set the value.
    $ws.Range($donorRef).Copy()
    $ws.Range($cellRef).PasteSpecial(-4122)  # xlPasteFormats
    $ws.Range($cellRef).Value = $val
}

function Set-Text0($cellRef) {
    # Re-use the text-cell style (copied from D14, which already shows "0").
    $ws.Range("D14").Copy()
    $ws.Range($cellRef).PasteSpecial(-4122)  # xlPasteFormats
    $ws.Range($cellRef).Formula = '="0"'
    $ws.Range($cellRef).Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)  # xlPasteValues -> literal shared string
}

function Set-TextStar($cellRef) {
    # Re-use the text-cell style (copied from E14, which already shows "***.*").
    $ws.Range("E14").Copy()
    $ws.Range($cellRef).PasteSpecial(-4122)  # xlPasteFormats
    $ws.Range($cellRef).Formula = '="***.*"'
    $ws.Range($cellRef).Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)  # xlPasteValues -> literal shared string
}

# --- Row 14 (Murder) --------------------------------------------------
Set-Text0 "C14"

# --- Row 15 (Rape) ------------------------------------------------------
Set-Num "C15" 1
Set-Num "F15" 5
Set-Text0 "G15"
Set-TextStar "H15"
Set-Num "I15" 14
Set-Num "K15" 27.272727272727
Set-Num "L15" -22.222222222222
Set-Num "M15" 7.692307692307
Set-Num "N15" -80.555555555555

# --- Row 16 (Robbery) ----------------------------------------------------
Set-Num "C16" 4
Set-Num "D16" 5
Set-Num "E16" -20
Set-Num "F16" 23
Set-Num "G16" 18
Set-Num "H16" 27.777777777777
Set-Num "I16" 109
Set-Num "J16" 87
Set-Num "K16" 25.287356321839
Set-Num "L16" 0.925925925925
Set-Num "M16" -47.846889952153
Set-Num "N16" -90.285204991087

# --- Row 17 (Felonious Assault) ------------------------------------------
Set-Num "C17" 7
Set-Num "D17" 5
Set-Num "E17" 40
Set-Num "F17" 35
Set-Num "G17" 21
Set-Num "H17" 66.666666666666
Set-Num "I17" 185
Set-Num "J17" 169
Set-Num "K17" 9.467455621301
Set-Num "L17" -7.035175879396
Set-Num "M17" -5.128205128205
Set-Num "N17" -61.855670103092

# --- Row 18 (Burglary) ----------------------------------------------------
Set-NumFrom "C18" "C17" 5
Set-Num "D18" 2
Set-Num "E18" 150
Set-Num "F18" 10
Set-Num "H18" 233.333333333333
Set-Num "I18" 69
Set-Num "J18" 71
Set-Num "K18" -2.816901408450
Set-Num "L18" -23.333333333333
Set-Num "M18" -52.413793103448
Set-Num "N18" -95.421366954213

# --- Row 19 (Grand Larceny) -----------------------------------------------
Set-Num "C19" 11
Set-Num "E19" 57.142857142857
Set-Num "F19" 42
Set-Num "G19" 32
Set-Num "H19" 31.25
Set-Num "I19" 282
Set-Num "J19" 295
Set-Num "K19" -4.406779661016
Set-Num "L19" 2.173913043478
Set-Num "M19" -14.024390243902
Set-Num "N19" -51.126516464471

# --- Row 20 (Grand Larceny Auto) ------------------------------------------
Set-Num "C20" 7
Set-Num "D20" 2
Set-Num "E20" 250
Set-Num "F20" 12
Set-Num "G20" 10
Set-Num "H20" 20
Set-Num "I20" 78
Set-Num "J20" 57
Set-Num "K20" 36.842105263157
Set-Num "L20" 13.043478260869
Set-Num "M20" -26.415094339622
Set-Num "N20" -94.372294372294

# --- Row 21 (TOTAL) --------------------------------------------------------
Set-Num "D21" 21
Set-Num "E21" 66.666666666666
Set-Num "F21" 128
Set-Num "G21" 84
Set-Num "H21" 52.380952380952
Set-Num "I21" 739
Set-Num "J21" 693
Set-Num "K21" 6.637806637806
Set-Num "L21" -2.890932982917
Set-Num "M21" -26.24750499002
Set-Num "N21" -85.711523588553

# --- Row 22 (Transit) -------------------------------------------------------
Set-Num "D22" 3
Set-Num "G22" 5
Set-Num "H22" -80
Set-Num "J22" 15
Set-Num "K22" -33.333333333333

# --- Row 23 (Housing) is unchanged this week --------------------------------

# --- Row 24 (Petit Larceny) -------------------------------------------------
Set-Num "C24" 40
Set-Num "D24" 50
Set-Num "E24" -20
Set-Num "F24" 133
Set-Num "G24" 154
Set-Num "H24" -13.636363636363
Set-Num "I24" 935
Set-Num "J24" 992
Set-Num "K24" -5.745967741935
Set-Num "L24" 35.507246376811
Set-Num "M24" 38.109305760709

# --- Row 25 (Retail Theft) --------------------------------------------------
Set-Num "C25" 23
Set-Num "D25" 24
Set-Num "E25" -4.166666666666
Set-Num "F25" 77
Set-Num "G25" 92
Set-Num "H25" -16.304347826087
Set-Num "I25" 575
Set-Num "J25" 574
Set-Num "K25" 0.174216027874
Set-Num "L25" 103.180212014134

# --- Row 26 (Misdemeanor Assault) -------------------------------------------
Set-Num "C26" 16
Set-Num "D26" 11
Set-Num "E26" 45.454545454545
Set-Num "F26" 55
Set-Num "G26" 41
Set-Num "H26" 34.146341463414
Set-Num "I26" 358
Set-Num "J26" 317
Set-Num "K26" 12.933753943217
Set-Num "L26" 11.875
Set-Num "M26" -12.682926829268

# --- Row 27 (UCR Rape*) ------------------------------------------------------
Set-Num "C27" 1
Set-Num "F27" 5
Set-Text0 "G27"
Set-TextStar "H27"
Set-Num "I27" 20
Set-Num "K27" 0
Set-Num "L27" -25.925925925925

# --- Row 28 (Other Sex Crimes) -----------------------------------------------
Set-Num "D28" 4
Set-Num "E28" -75
Set-Num "F28" 4
Set-Num "G28" 8
Set-Num "H28" -50
Set-Num "I28" 39
Set-Num "J28" 35
Set-Num "K28" 11.428571428571
Set-Num "L28" 11.428571428571

# --- Row 29 (Shooting Victims) ------------------------------------------------
Set-Num "C29" 2
Set-Num "F29" 4
Set-Text0 "G29"
Set-TextStar "H29"
Set-Num "I29" 7
Set-Num "K29" 0
Set-Num "L29" 0
Set-Num "M29" -70.833333333333
Set-Num "N29" -91.358024691358

# --- Row 30 (Shooting Incidents) ----------------------------------------------
Set-Num "F30" 3
Set-Text0 "G30"
Set-TextStar "H30"
Set-Num "I30" 6
Set-Num "K30" 0
Set-Num "L30" -14.285714285714
Set-Num "M30" -68.421052631578
Set-Num "N30" -90.322580645161

# --- Row 31 (Hate Crimes) ------------------------------------------------------
Set-Text0 "F31"
Set-Num "L31" -10
